$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new parent sample value (containing a comma) to the next empty row in column A
$ws.Range("A6").Value = '"parent,06"'

# Move/set the active selection as it appears after the edit in the saved file
$ws.Range("C8").Select()
